$d = $word.ActiveDocument

$replacements = @(
    @("40÷6=", "20÷6="),
    @("38÷3=", "19÷8="),
    @("98÷9=", "32÷3="),
    @("68÷3=", "31÷5="),
    @("97÷8=", "19÷8="),
    @("15÷3=", "34÷9="),
    @("65÷5=", "77÷8="),
    @("21÷9=", "83÷6="),
    @("23÷3=", "11÷2="),
    @("54÷9=", "16÷2="),
    @("74÷3=", "44÷7="),
    @("98÷2=", "44÷7="),
    @("31÷2=", "24÷6="),
    @("94÷6=", "81÷7="),
    @("33÷2=", "31÷4="),
    @("50÷2=", "51÷4="),
    @("29÷5=", "65÷4="),
    @("25÷2=", "44÷8="),
    @("64÷9=", "37÷5="),
    @("96÷5=", "17÷8="),
    @("18÷6=", "71÷3="),
    @("22÷7=", "50÷2="),
    @("68÷5=", "86÷5="),
    @("32÷4=", "68÷8="),
    @("92÷2=", "94÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
